$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 8671.799999999999
$ws.Range("J4").Value = 15999.5
$ws.Range("L4").Value = 15999.5
$ws.Range("N4").Value = -16227.5
# Row 116
$ws.Range("H116").Value = 3707.5557
$ws.Range("I116").Value = 2971.7144
$ws.Range("J116").Value = 4500
$ws.Range("K116").Value = 2971.7144
$ws.Range("L116").Value = 4500
$ws.Range("M116").Value = 470.2856000000002
$ws.Range("N116").Value = -11384
# Row 138
$ws.Range("H138").Value = 3926.1592
$ws.Range("I138").Value = 2858.0908
$ws.Range("J138").Value = 4282.1816
$ws.Range("K138").Value = 8574.2724
$ws.Range("L138").Value = 12846.5448
$ws.Range("M138").Value = -3434.2724
$ws.Range("N138").Value = -23126.5448

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 8
$ws.Range("H8").Value = 2247
$ws.Range("I8").Value = 2247
$ws.Range("K8").Value = 2247
$ws.Range("M8").Value = -2103
# Row 16
$ws.Range("H16").Value = 278
$ws.Range("I16").Value = 278
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 278
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 9
$ws.Range("N16").ClearContents()
# Row 32
$ws.Range("H32").Value = 3352604.8
$ws.Range("I32").Value = 3186814.2
$ws.Range("K32").Value = 3186814.2
$ws.Range("M32").Value = -3186527.2
# Row 61
$ws.Range("H61").Value = 5362.52
$ws.Range("I61").Value = 5433.174
$ws.Range("K61").Value = 5433.174
$ws.Range("M61").Value = -5221.174
# Row 74
$ws.Range("H74").Value = 1056.2858
$ws.Range("I74").Value = 978.8
$ws.Range("J74").Value = 1250
$ws.Range("K74").Value = 978.8
$ws.Range("L74").Value = 1250
$ws.Range("M74").Value = -104.8
$ws.Range("N74").Value = -2998
# Row 77
$ws.Range("H77").Value = 1056.2858
$ws.Range("I77").Value = 978.8
$ws.Range("J77").Value = 1250
$ws.Range("K77").Value = 4894
$ws.Range("L77").Value = 6250
$ws.Range("M77").Value = -526
$ws.Range("N77").Value = -14986
# Row 102
$ws.Range("H102").Value = 606.1429000000001
$ws.Range("I102").Value = 548.8
$ws.Range("J102").Value = 749.5
$ws.Range("K102").Value = 548.8
$ws.Range("L102").Value = 749.5
$ws.Range("M102").Value = 1073.2
$ws.Range("N102").Value = -3993.5
# Row 122
$ws.Range("H122").Value = 18020.895
$ws.Range("I122").Value = 18020.895
$ws.Range("K122").Value = 54062.685
$ws.Range("M122").Value = -51612.685
# Row 132
$ws.Range("H132").Value = 2190.1
$ws.Range("I132").Value = 2186.4285
$ws.Range("K132").Value = 6559.2855
$ws.Range("M132").Value = -4029.2855
# Row 136
$ws.Range("H136").Value = 5362.52
$ws.Range("I136").Value = 5433.174
$ws.Range("K136").Value = 16299.522
$ws.Range("M136").Value = -13749.522

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 305.5
$ws.Range("I5").Value = 99.5
$ws.Range("K5").Value = 99.5
$ws.Range("M5").Value = 13.5
# Row 57
$ws.Range("H57").Value = 89998.5
$ws.Range("I57").Value = 89998.5
$ws.Range("K57").Value = 89998.5
$ws.Range("M57").Value = -89278.5
# Row 105
$ws.Range("H105").Value = 2408.5715
$ws.Range("I105").Value = 2241.5
$ws.Range("K105").Value = 2241.5
$ws.Range("M105").Value = -494.5
# Row 134
$ws.Range("H134").Value = 4159.0527
$ws.Range("I134").Value = 4334.5557
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 13003.6671
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -10468.6671
$ws.Range("N134").Value = -8070
# Row 136
$ws.Range("H136").Value = 89998.5
$ws.Range("I136").Value = 89998.5
$ws.Range("K136").Value = 89998.5
$ws.Range("M136").Value = -84898.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 171.90909
$ws.Range("I7").Value = 89.09999999999999
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 89.09999999999999
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = 23.90000000000001
$ws.Range("N7").Value = -1226
# Row 22
$ws.Range("H22").Value = 2284.2856
$ws.Range("I22").Value = 2498.75
$ws.Range("J22").Value = 1998.3334
$ws.Range("K22").Value = 2498.75
$ws.Range("L22").Value = 1998.3334
$ws.Range("M22").Value = -2148.75
$ws.Range("N22").Value = -2698.3334
# Row 31
$ws.Range("H31").Value = 1360.2273
$ws.Range("J31").Value = 1278.7222
$ws.Range("L31").Value = 1278.7222
$ws.Range("N31").Value = -1868.7222
# Row 34
$ws.Range("H34").Value = 1360.2273
$ws.Range("J34").Value = 1278.7222
$ws.Range("L34").Value = 1278.7222
$ws.Range("N34").Value = -1682.7222
# Row 94
$ws.Range("H94").Value = 96307.75
$ws.Range("J94").Value = 4607.1113
$ws.Range("L94").Value = 4607.1113
$ws.Range("N94").Value = -5509.1113
# Row 132
$ws.Range("H132").Value = 7990.6875
$ws.Range("I132").Value = 8000.067
$ws.Range("J132").Value = 7850
$ws.Range("K132").Value = 24000.201
$ws.Range("L132").Value = 23550
$ws.Range("M132").Value = -21470.201
$ws.Range("N132").Value = -28610

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 4314.1763
$ws.Range("J68").Value = 4314.1763
$ws.Range("L68").Value = 12942.5289
$ws.Range("N68").Value = -14564.5289
# Row 71
$ws.Range("H71").Value = 4314.1763
$ws.Range("J71").Value = 4314.1763
$ws.Range("L71").Value = 38827.5867
$ws.Range("N71").Value = -46939.5867
# Row 88
$ws.Range("H88").Value = 12999.5
$ws.Range("J88").Value = 12999.5
$ws.Range("L88").Value = 38998.5
$ws.Range("N88").Value = -39854.5
# Row 91
$ws.Range("H91").Value = 12999.5
$ws.Range("J91").Value = 12999.5
$ws.Range("L91").Value = 38998.5
$ws.Range("N91").Value = -41962.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 917258.2
$ws.Range("I3").Value = 1100660
$ws.Range("K3").Value = 1100660
$ws.Range("M3").Value = -1100544
# Row 14
$ws.Range("H14").Value = 62840.562
$ws.Range("I14").Value = 71774.92999999999
$ws.Range("K14").Value = 71774.92999999999
$ws.Range("M14").Value = -71606.92999999999
# Row 21
$ws.Range("H21").Value = 6005
$ws.Range("I21").Value = 6005
$ws.Range("K21").Value = 6005
$ws.Range("M21").Value = -5832
# Row 30
$ws.Range("H30").Value = 6005
$ws.Range("I30").Value = 6005
$ws.Range("K30").Value = 6005
$ws.Range("M30").Value = -5900
# Row 102
$ws.Range("H102").Value = 2624.4167
$ws.Range("I102").Value = 2247.2856
$ws.Range("K102").Value = 2247.2856
$ws.Range("M102").Value = -625.2856000000002
# Row 132
$ws.Range("H132").Value = 3591.5
$ws.Range("I132").Value = 3591.5
$ws.Range("K132").Value = 10774.5
$ws.Range("M132").Value = -8244.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 2194.818
$ws.Range("I55").Value = 2491.4
$ws.Range("K55").Value = 2491.4
$ws.Range("M55").Value = -2318.4
# Row 82
$ws.Range("H82").Value = 917.5
$ws.Range("I82").Value = 898.5454999999999
$ws.Range("J82").Value = 959.2
$ws.Range("K82").Value = 898.5454999999999
$ws.Range("L82").Value = 959.2
$ws.Range("M82").Value = -537.5454999999999
$ws.Range("N82").Value = -1681.2
# Row 85
$ws.Range("H85").Value = 917.5
$ws.Range("I85").Value = 898.5454999999999
$ws.Range("J85").Value = 959.2
$ws.Range("K85").Value = 898.5454999999999
$ws.Range("L85").Value = 959.2
$ws.Range("M85").Value = 349.4545000000001
$ws.Range("N85").Value = -3455.2
# Row 93
$ws.Range("H93").Value = 2900
$ws.Range("I93").Value = 2900
$ws.Range("K93").Value = 2900
$ws.Range("M93").Value = -1652

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 7623.25
# Row 65
$ws.Range("H65").Value = 7623.25
# Row 100
$ws.Range("H100").Value = 100000000
$ws.Range("I100").Value = 100000000
$ws.Range("K100").Value = 200000000
$ws.Range("M100").Value = -199999459
# Row 132
$ws.Range("H132").Value = 1939
$ws.Range("I132").Value = 2001.7142
$ws.Range("K132").Value = 6005.142599999999
$ws.Range("M132").Value = -3475.142599999999
